$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1063.5
$ws.Range("I32").Value = 1450
$ws.Range("J32").Value = 677
$ws.Range("K32").Value = 1450
$ws.Range("L32").Value = 677
$ws.Range("M32").Value = -1124
$ws.Range("N32").Value = -1329
$ws.Range("H33").Value = 149.41176
$ws.Range("I33").Value = 217.42857
$ws.Range("J33").Value = 101.8
$ws.Range("K33").Value = 217.42857
$ws.Range("L33").Value = 101.8
$ws.Range("M33").Value = 11.57142999999999
$ws.Range("N33").Value = -559.8
$ws.Range("H76").Value = 4793.6
$ws.Range("H79").Value = 4793.6
$ws.Range("H111").Value = 1724.875
$ws.Range("I111").Value = 1602.4
$ws.Range("J111").Value = 1929
$ws.Range("K111").Value = 4807.200000000001
$ws.Range("L111").Value = 5787
$ws.Range("M111").Value = -1740.200000000001
$ws.Range("N111").Value = -11921
$ws.Range("H112").Value = 3128.6956
$ws.Range("J112").Value = 3128.6956
$ws.Range("L112").Value = 9386.086800000001
$ws.Range("N112").Value = -11602.0868
$ws.Range("H137").Value = 1737.5
$ws.Range("I137").Value = 1280
$ws.Range("K137").Value = 3840
$ws.Range("M137").Value = -1290
$ws.Range("H138").Value = 5014.2856
$ws.Range("I138").Value = 3899.75
$ws.Range("K138").Value = 11699.25
$ws.Range("M138").Value = -6559.25

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3869.3333
$ws.Range("I2").Value = 798
$ws.Range("K2").Value = 798
$ws.Range("M2").Value = -685
$ws.Range("H4").Value = 41
$ws.Range("I4").Value = 41
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 41
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 75
$ws.Range("N4").ClearContents()
$ws.Range("H22").Value = 4000
$ws.Range("I22").Value = 3500
$ws.Range("K22").Value = 3500
$ws.Range("M22").Value = -3201
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("N25").ClearContents()
$ws.Range("H32").Value = 2344.5908
$ws.Range("I32").Value = 2344.5908
$ws.Range("K32").Value = 2344.5908
$ws.Range("M32").Value = -2057.5908
$ws.Range("H45").Value = 4799.8
$ws.Range("I45").Value = 1499.5
$ws.Range("J45").Value = 7000
$ws.Range("K45").Value = 1499.5
$ws.Range("L45").Value = 7000
$ws.Range("M45").Value = -1122.5
$ws.Range("N45").Value = -7754
$ws.Range("H63").Value = 2099.8333
$ws.Range("H66").Value = 2099.8333
$ws.Range("H95").Value = 2900
$ws.Range("J95").Value = 2900
$ws.Range("L95").Value = 2900
$ws.Range("N95").Value = -8392
$ws.Range("H97").Value = 90.42856999999999
$ws.Range("I97").Value = 95
$ws.Range("K97").Value = 95
$ws.Range("M97").Value = 401
$ws.Range("H116").Value = 3869.3333
$ws.Range("I116").Value = 798
$ws.Range("K116").Value = 798
$ws.Range("M116").Value = 1496

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3869.3333
$ws.Range("I3").Value = 798
$ws.Range("K3").Value = 798
$ws.Range("M3").Value = -684
$ws.Range("H94").Value = 2996
$ws.Range("I94").Value = 491
$ws.Range("K94").Value = 491
$ws.Range("M94").Value = -40
$ws.Range("H99").Value = 2350
$ws.Range("I99").Value = 2300
$ws.Range("J99").Value = 2400
$ws.Range("K99").Value = 2300
$ws.Range("L99").Value = 2400
$ws.Range("M99").Value = -802
$ws.Range("N99").Value = -5396
$ws.Range("H134").Value = 4573.8335
$ws.Range("I134").Value = 4339.6
$ws.Range("J134").Value = 4741.143
$ws.Range("K134").Value = 13018.8
$ws.Range("L134").Value = 14223.429
$ws.Range("M134").Value = -10483.8
$ws.Range("N134").Value = -19293.429

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2040
$ws.Range("I31").Value = 2040
$ws.Range("K31").Value = 2040
$ws.Range("M31").Value = -1745
$ws.Range("H34").Value = 2040
$ws.Range("I34").Value = 2040
$ws.Range("K34").Value = 2040
$ws.Range("M34").Value = -1838
$ws.Range("H86").Value = 9960499
$ws.Range("I86").Value = 9960499
$ws.Range("K86").Value = 9960499
$ws.Range("M86").Value = -9959376
$ws.Range("H89").Value = 9960499
$ws.Range("I89").Value = 9960499
$ws.Range("K89").Value = 49802495
$ws.Range("M89").Value = -49796879
$ws.Range("H112").Value = 51499.5
$ws.Range("J112").Value = 51499.5
$ws.Range("L112").Value = 51499.5
$ws.Range("N112").Value = -54453.5
$ws.Range("H122").Value = 2245
$ws.Range("I122").Value = 1860.5
$ws.Range("J122").Value = 3014
$ws.Range("K122").Value = 5581.5
$ws.Range("L122").Value = 9042
$ws.Range("M122").Value = -3131.5
$ws.Range("N122").Value = -13942
$ws.Range("H141").Value = 400564
$ws.Range("J141").Value = 400564
$ws.Range("L141").Value = 400564
$ws.Range("N141").Value = -410924

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 95
$ws.Range("I11").Value = 95
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 285
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -145
$ws.Range("N11").ClearContents()
$ws.Range("H26").Value = 28.333334
$ws.Range("I26").Value = 28.333334
$ws.Range("K26").Value = 85.00000199999999
$ws.Range("M26").Value = 202.999998
$ws.Range("H81").Value = 10066.667
$ws.Range("I81").Value = 250
$ws.Range("J81").Value = 14975
$ws.Range("K81").Value = 750
$ws.Range("L81").Value = 44925
$ws.Range("M81").Value = 373
$ws.Range("N81").Value = -47171
$ws.Range("H84").Value = 10066.667
$ws.Range("I84").Value = 250
$ws.Range("J84").Value = 14975
$ws.Range("K84").Value = 2250
$ws.Range("L84").Value = 134775
$ws.Range("M84").Value = 3366
$ws.Range("N84").Value = -146007

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9177.799999999999
$ws.Range("I70").Value = 9972.25
$ws.Range("K70").Value = 9972.25
$ws.Range("M70").Value = -9702.25
$ws.Range("H73").Value = 9177.799999999999
$ws.Range("I73").Value = 9972.25
$ws.Range("K73").Value = 9972.25
$ws.Range("M73").Value = -9036.25

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H22").Value = 1416.1428
$ws.Range("J22").Value = 1151.75
$ws.Range("L22").Value = 1151.75
$ws.Range("N22").Value = -1741.75
$ws.Range("H27").Value = 1416.1428
$ws.Range("J27").Value = 1151.75
$ws.Range("L27").Value = 1151.75
$ws.Range("N27").Value = -1365.75
$ws.Range("H40").Value = 2081.25
$ws.Range("I40").Value = 2081.25
$ws.Range("K40").Value = 2081.25
$ws.Range("M40").Value = -1945.25
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H100").Value = 2579
$ws.Range("I100").Value = 2579
$ws.Range("K100").Value = 2579
$ws.Range("M100").Value = -2038

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 11616936
$ws.Range("I100").Value = 23233022
$ws.Range("J100").Value = 849.6667
$ws.Range("K100").Value = 46466044
$ws.Range("L100").Value = 1699.3334
$ws.Range("M100").Value = -46465503
$ws.Range("N100").Value = -2781.3334
$ws.Range("H112").Value = 50188.5
$ws.Range("J112").Value = 50188.5
$ws.Range("L112").Value = 50188.5
$ws.Range("N112").Value = -53142.5
$ws.Range("H113").Value = 1307.875
$ws.Range("I113").Value = 977.3333
$ws.Range("K113").Value = 2931.9999
$ws.Range("M113").Value = -761.9998999999998
$ws.Range("H126").Value = 897.8461
$ws.Range("I126").Value = 933.9091
$ws.Range("J126").Value = 699.5
$ws.Range("K126").Value = 2801.7273
$ws.Range("L126").Value = 2098.5
$ws.Range("M126").Value = -331.7273
$ws.Range("N126").Value = -7038.5
$ws.Range("H136").Value = 959.6818
$ws.Range("I136").Value = 965.7
$ws.Range("K136").Value = 2897.1
$ws.Range("M136").Value = -347.1000000000004
